$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a new data row (row 68) below the existing dataset ---
# Match the existing date/time number format used by the rest of the
# column (style index carrying numFmt "dd/mm/yy hh:mm") so the new
# cells land in the same style as the rows above them instead of
# triggering a brand-new auto-detected number format.
$dtFormat = "dd/mm/yy\ hh:mm"
$newRow = 68

$ws.Range("A$newRow`:F$newRow").NumberFormat = $dtFormat

$ws.Range("A$newRow").Value = 42474.375
$ws.Range("B$newRow").Value = 42474.375
$ws.Range("C$newRow").Value = 42474.895833333336
$ws.Range("D$newRow").Value = 42474.895833333336
$ws.Range("E$newRow").Value = 42474.895833333336
$ws.Range("F$newRow").Value = 42474.895833333336

# --- Widen columns B and C (they stop being auto "best fit" and become
# explicit, wider, user-set widths) ---
$ws.Columns("B").ColumnWidth = 18.85546875
$ws.Columns("C").ColumnWidth = 22.42578125

# --- Move the selection to the empty row right after the newly added
# data, matching where the user's cursor ended up after entering the
# row (select the whole next row) ---
$ws.Range("A69:XFD69").Select()
